$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.286709999999999
$ws.Range("H2").Value = 21.86013
$ws.Range("I2").Value = 0.1632739668438106
$ws.Range("J2").Value = 0.1632739668438107
$ws.Range("M2").Value = 50.26671733333333
$ws.Range("N2").Value = 150.800152
$ws.Range("O2").Value = 0.2619764206727233
$ws.Range("P2").Value = 0.2619764206727234
$ws.Range("Q2").Value = 366.2789918599733
$ws.Range("R2").Value = 3296.51092673976
$ws.Range("S2").Value = 0.04277392942277841
$ws.Range("T2").Value = 0.04277392942277843
$ws.Range("G3").Value = 7.286709999999999
$ws.Range("H3").Value = 21.86013
$ws.Range("I3").Value = 0.1632739668438106
$ws.Range("J3").Value = 0.1632739668438107
$ws.Range("O3").Value = 0.01622594841727
$ws.Range("P3").Value = 0.01622594841727001
$ws.Range("Q3").Value = 22.68610286753333
$ws.Range("R3").Value = 204.1749258078
$ws.Range("S3").Value = 0.002649274963890724
$ws.Range("T3").Value = 0.002649274963890725
$ws.Range("G4").Value = 7.286709999999999
$ws.Range("H4").Value = 21.86013
$ws.Range("I4").Value = 0.1632739668438106
$ws.Range("J4").Value = 0.1632739668438107
$ws.Range("M4").Value = 80.63290666666667
$ws.Range("N4").Value = 241.89872
$ws.Range("O4").Value = 0.4202367172077739
$ws.Range("P4").Value = 0.4202367172077739
$ws.Range("Q4").Value = 587.5486073370666
$ws.Range("R4").Value = 5287.9374660336
$ws.Range("S4").Value = 0.0686137158319339
$ws.Range("T4").Value = 0.0686137158319339
$ws.Range("G5").Value = 7.286709999999999
$ws.Range("H5").Value = 21.86013
$ws.Range("I5").Value = 0.1632739668438106
$ws.Range("J5").Value = 0.1632739668438107
$ws.Range("M5").Value = 2.274154
$ws.Range("N5").Value = 6.822462
$ws.Range("O5").Value = 0.01185227038057408
$ws.Range("P5").Value = 0.01185227038057408
$ws.Range("Q5").Value = 16.57110069334
$ws.Range("R5").Value = 149.13990624006
$ws.Range("S5").Value = 0.001935167201141731
$ws.Range("T5").Value = 0.001935167201141732
$ws.Range("G6").Value = 7.286709999999999
$ws.Range("H6").Value = 21.86013
$ws.Range("I6").Value = 0.1632739668438106
$ws.Range("J6").Value = 0.1632739668438107
$ws.Range("M6").Value = 55.58783666666667
$ws.Range("N6").Value = 166.76351
$ws.Range("O6").Value = 0.2897086433216586
$ws.Range("P6").Value = 0.2897086433216586
$ws.Range("Q6").Value = 405.0524453173666
$ws.Range("R6").Value = 3645.4720078563
$ws.Range("S6").Value = 0.04730187942406585
$ws.Range("T6").Value = 0.04730187942406586
$ws.Range("I7").Value = 0.05961564640488622
$ws.Range("J7").Value = 0.05961564640488623
$ws.Range("M7").Value = 50.26671733333333
$ws.Range("N7").Value = 150.800152
$ws.Range("O7").Value = 0.2619764206727233
$ws.Range("P7").Value = 0.2619764206727234
$ws.Range("Q7").Value = 133.7381536466916
$ws.Range("R7").Value = 1203.643382820224
$ws.Range("S7").Value = 0.0156178936612428
$ws.Range("T7").Value = 0.0156178936612428
$ws.Range("I8").Value = 0.05961564640488622
$ws.Range("J8").Value = 0.05961564640488623
$ws.Range("O8").Value = 0.01622594841727
$ws.Range("P8").Value = 0.01622594841727001
$ws.Range("S8").Value = 0.0009673204034278918
$ws.Range("T8").Value = 0.0009673204034278921
$ws.Range("I9").Value = 0.05961564640488622
$ws.Range("J9").Value = 0.05961564640488623
$ws.Range("M9").Value = 80.63290666666667
$ws.Range("N9").Value = 241.89872
$ws.Range("O9").Value = 0.4202367172077739
$ws.Range("P9").Value = 0.4202367172077739
$ws.Range("Q9").Value = 214.5295462454044
$ws.Range("R9").Value = 1930.76591620864
$ws.Range("S9").Value = 0.02505268353940882
$ws.Range("T9").Value = 0.02505268353940882
$ws.Range("I10").Value = 0.05961564640488622
$ws.Range("J10").Value = 0.05961564640488623
$ws.Range("M10").Value = 2.274154
$ws.Range("N10").Value = 6.822462
$ws.Range("O10").Value = 0.01185227038057408
$ws.Range("P10").Value = 0.01185227038057408
$ws.Range("Q10").Value = 6.050547423882666
$ws.Range("R10").Value = 54.454926814944
$ws.Range("S10").Value = 0.0007065807601034107
$ws.Range("T10").Value = 0.0007065807601034108
$ws.Range("I11").Value = 0.05961564640488622
$ws.Range("J11").Value = 0.05961564640488623
$ws.Range("M11").Value = 55.58783666666667
$ws.Range("N11").Value = 166.76351
$ws.Range("O11").Value = 0.2897086433216586
$ws.Range("P11").Value = 0.2897086433216586
$ws.Range("Q11").Value = 147.8953676587911
$ws.Range("R11").Value = 1331.05830892912
$ws.Range("S11").Value = 0.01727116804070331
$ws.Range("T11").Value = 0.01727116804070331
$ws.Range("G12").Value = 19.00851733333333
$ws.Range("H12").Value = 57.025552
$ws.Range("I12").Value = 0.4259255588369328
$ws.Range("J12").Value = 0.4259255588369329
$ws.Range("M12").Value = 50.26671733333333
$ws.Range("N12").Value = 150.800152
$ws.Range("O12").Value = 0.2619764206727233
$ws.Range("P12").Value = 0.2619764206727234
$ws.Range("Q12").Value = 955.4957677204338
$ws.Range("R12").Value = 8599.461909483905
$ws.Range("S12").Value = 0.1115824533771291
$ws.Range("T12").Value = 0.1115824533771291
$ws.Range("G13").Value = 19.00851733333333
$ws.Range("H13").Value = 57.025552
$ws.Range("I13").Value = 0.4259255588369328
$ws.Range("J13").Value = 0.4259255588369329
$ws.Range("O13").Value = 0.01622594841727
$ws.Range("P13").Value = 0.01622594841727001
$ws.Range("Q13").Value = 59.18023080145777
$ws.Range("R13").Value = 532.62207721312
$ws.Range("S13").Value = 0.006911046147284971
$ws.Range("T13").Value = 0.006911046147284975
$ws.Range("G14").Value = 19.00851733333333
$ws.Range("H14").Value = 57.025552
$ws.Range("I14").Value = 0.4259255588369328
$ws.Range("J14").Value = 0.4259255588369329
$ws.Range("M14").Value = 80.63290666666667
$ws.Range("N14").Value = 241.89872
$ws.Range("O14").Value = 0.4202367172077739
$ws.Range("P14").Value = 0.4202367172077739
$ws.Range("Q14").Value = 1532.712004010382
$ws.Range("R14").Value = 13794.40803609344
$ws.Range("S14").Value = 0.1789895586205192
$ws.Range("T14").Value = 0.1789895586205192
$ws.Range("G15").Value = 19.00851733333333
$ws.Range("H15").Value = 57.025552
$ws.Range("I15").Value = 0.4259255588369328
$ws.Range("J15").Value = 0.4259255588369329
$ws.Range("M15").Value = 2.274154
$ws.Range("N15").Value = 6.822462
$ws.Range("O15").Value = 0.01185227038057408
$ws.Range("P15").Value = 0.01185227038057408
$ws.Range("Q15").Value = 43.22829572766933
$ws.Range("R15").Value = 389.054661549024
$ws.Range("S15").Value = 0.005048184885332442
$ws.Range("T15").Value = 0.005048184885332443
$ws.Range("G16").Value = 19.00851733333333
$ws.Range("H16").Value = 57.025552
$ws.Range("I16").Value = 0.4259255588369328
$ws.Range("J16").Value = 0.4259255588369329
$ws.Range("M16").Value = 55.58783666666667
$ws.Range("N16").Value = 166.76351
$ws.Range("O16").Value = 0.2897086433216586
$ws.Range("P16").Value = 0.2897086433216586
$ws.Range("Q16").Value = 1056.642356800836
$ws.Range("R16").Value = 9509.781211207521
$ws.Range("S16").Value = 0.1233943158066671
$ws.Range("T16").Value = 0.1233943158066671
$ws.Range("G17").Value = 0.1146546666666667
$ws.Range("H17").Value = 0.343964
$ws.Range("I17").Value = 0.002569077436020028
$ws.Range("J17").Value = 0.002569077436020028
$ws.Range("M17").Value = 50.26671733333333
$ws.Range("N17").Value = 150.800152
$ws.Range("O17").Value = 0.2619764206727233
$ws.Range("P17").Value = 0.2619764206727234
$ws.Range("Q17").Value = 5.763313720280889
$ws.Range("R17").Value = 51.869823482528
$ws.Range("S17").Value = 0.0006730377111195841
$ws.Range("T17").Value = 0.0006730377111195844
$ws.Range("G18").Value = 0.1146546666666667
$ws.Range("H18").Value = 0.343964
$ws.Range("I18").Value = 0.002569077436020028
$ws.Range("J18").Value = 0.002569077436020028
$ws.Range("O18").Value = 0.01622594841727
$ws.Range("P18").Value = 0.01622594841727001
$ws.Range("Q18").Value = 0.3569604886488889
$ws.Range("R18").Value = 3.21264439784
$ws.Range("S18").Value = 0.00004168571795683324
$ws.Range("T18").Value = 0.00004168571795683326
$ws.Range("G19").Value = 0.1146546666666667
$ws.Range("H19").Value = 0.343964
$ws.Range("I19").Value = 0.002569077436020028
$ws.Range("J19").Value = 0.002569077436020028
$ws.Range("M19").Value = 80.63290666666667
$ws.Range("N19").Value = 241.89872
$ws.Range("O19").Value = 0.4202367172077739
$ws.Range("P19").Value = 0.4202367172077739
$ws.Range("Q19").Value = 9.244939036231111
$ws.Range("R19").Value = 83.20445132607999
$ws.Range("S19").Value = 0.001079620667965621
$ws.Range("T19").Value = 0.001079620667965621
$ws.Range("G20").Value = 0.1146546666666667
$ws.Range("H20").Value = 0.343964
$ws.Range("I20").Value = 0.002569077436020028
$ws.Range("J20").Value = 0.002569077436020028
$ws.Range("M20").Value = 2.274154
$ws.Range("N20").Value = 6.822462
$ws.Range("O20").Value = 0.01185227038057408
$ws.Range("P20").Value = 0.01185227038057408
$ws.Range("Q20").Value = 0.2607423688186666
$ws.Range("R20").Value = 2.346681319368
$ws.Range("S20").Value = 0.00003044940040034138
$ws.Range("T20").Value = 0.00003044940040034138
$ws.Range("G21").Value = 0.1146546666666667
$ws.Range("H21").Value = 0.343964
$ws.Range("I21").Value = 0.002569077436020028
$ws.Range("J21").Value = 0.002569077436020028
$ws.Range("M21").Value = 55.58783666666667
$ws.Range("N21").Value = 166.76351
$ws.Range("O21").Value = 0.2897086433216586
$ws.Range("P21").Value = 0.2897086433216586
$ws.Range("Q21").Value = 6.373404883737778
$ws.Range("R21").Value = 57.36064395364
$ws.Range("S21").Value = 0.0007442839385776474
$ws.Range("T21").Value = 0.0007442839385776474
$ws.Range("G22").Value = 15.55827866666667
$ws.Range("H22").Value = 46.674836
$ws.Range("I22").Value = 0.3486157504783503
$ws.Range("J22").Value = 0.3486157504783503
$ws.Range("M22").Value = 50.26671733333333
$ws.Range("N22").Value = 150.800152
$ws.Range("O22").Value = 0.2619764206727233
$ws.Range("P22").Value = 0.2619764206727234
$ws.Range("Q22").Value = 782.0635959305636
$ws.Range("R22").Value = 7038.572363375071
$ws.Range("S22").Value = 0.09132910650045342
$ws.Range("T22").Value = 0.09132910650045346
$ws.Range("G23").Value = 15.55827866666667
$ws.Range("H23").Value = 46.674836
$ws.Range("I23").Value = 0.3486157504783503
$ws.Range("J23").Value = 0.3486157504783503
$ws.Range("O23").Value = 0.01622594841727
$ws.Range("P23").Value = 0.01622594841727001
$ws.Range("Q23").Value = 48.43841874779555
$ws.Range("R23").Value = 435.94576873016
$ws.Range("S23").Value = 0.005656621184709582
$ws.Range("T23").Value = 0.005656621184709584
$ws.Range("G24").Value = 15.55827866666667
$ws.Range("H24").Value = 46.674836
$ws.Range("I24").Value = 0.3486157504783503
$ws.Range("J24").Value = 0.3486157504783503
$ws.Range("M24").Value = 80.63290666666667
$ws.Range("N24").Value = 241.89872
$ws.Range("O24").Value = 0.4202367172077739
$ws.Range("P24").Value = 0.4202367172077739
$ws.Range("Q24").Value = 1254.509231623324
$ws.Range("R24").Value = 11290.58308460992
$ws.Range("S24").Value = 0.1465011385479463
$ws.Range("T24").Value = 0.1465011385479464
$ws.Range("G25").Value = 15.55827866666667
$ws.Range("H25").Value = 46.674836
$ws.Range("I25").Value = 0.3486157504783503
$ws.Range("J25").Value = 0.3486157504783503
$ws.Range("M25").Value = 2.274154
$ws.Range("N25").Value = 6.822462
$ws.Range("O25").Value = 0.01185227038057408
$ws.Range("P25").Value = 0.01185227038057408
$ws.Range("Q25").Value = 35.38192166291466
$ws.Range("R25").Value = 318.437294966232
$ws.Range("S25").Value = 0.004131888133596155
$ws.Range("T25").Value = 0.004131888133596156
$ws.Range("G26").Value = 15.55827866666667
$ws.Range("H26").Value = 46.674836
$ws.Range("I26").Value = 0.3486157504783503
$ws.Range("J26").Value = 0.3486157504783503
$ws.Range("M26").Value = 55.58783666666667
$ws.Range("N26").Value = 166.76351
$ws.Range("O26").Value = 0.2897086433216586
$ws.Range("P26").Value = 0.2897086433216586
$ws.Range("Q26").Value = 864.8510533371511
$ws.Range("R26").Value = 7783.65948003436
$ws.Range("S26").Value = 0.1009969961116447
$ws.Range("T26").Value = 0.1009969961116447

Write-Host "Updated 278 cells with new TPM-derived values"
